$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (Strike#) values for rows 2-19, column G
$kValues = @{
    2  = 2
    3  = 6
    4  = 4
    5  = 2
    6  = 4
    7  = 2
    8  = 3
    9  = 1
    10 = 1
    11 = 2
    12 = 8
    13 = 7
    14 = 3
    15 = 2
    16 = 1
    17 = 2
    18 = 4
    19 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
